$d = $word.ActiveDocument

$replacements = @(
    @("2024-04-08 Monday", "2024-04-09 Tuesday"),
    @("41×43=", "26×50="),
    @("43×47=", "20×77="),
    @("85×80=", "14×95="),
    @("66×21=", "17×65="),
    @("75×97=", "50×29="),
    @("42×82=", "35×13="),
    @("31×57=", "25×37="),
    @("26×39=", "58×82="),
    @("33×36=", "49×63="),
    @("86×92=", "44×67="),
    @("38×26=", "48×39="),
    @("17×91=", "33×37="),
    @("22×97=", "86×45="),
    @("21×52=", "79×63="),
    @("18×17=", "98×42="),
    @("99×44=", "73×61="),
    @("40×51=", "87×21="),
    @("13×68=", "60×43="),
    @("28×39=", "68×39="),
    @("62×40=", "78×79="),
    @("59×56=", "80×19="),
    @("91×53=", "94×21="),
    @("30×80=", "71×65="),
    @("42×70=", "89×77="),
    @("49×12=", "73×42=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
